$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    4258.591237880026,
    4234.383934797102,
    4234.383934797102,
    4202.57169580428,
    3997.946913772582,
    3997.946913772582,
    3997.946913772582,
    3997.946913772582,
    3991.189763750467,
    3923.511757890763,
    3923.511757890763
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
